# chore: update Sheets via scheduled runner
# Refreshes market-price derived figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# for a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 132.76471
$ws.Range("I39").Value = 33.1
$ws.Range("J39").Value = 275.14285
$ws.Range("K39").Value = 99.30000000000001
$ws.Range("L39").Value = 825.4285500000001
$ws.Range("M39").Value = 196.7
$ws.Range("N39").Value = -1417.42855
$ws.Range("H64").Value = 9999.4
$ws.Range("J64").Value = 9999.4
$ws.Range("L64").Value = 9999.4
$ws.Range("N64").Value = -10495.4
$ws.Range("H67").Value = 9999.4
$ws.Range("J67").Value = 9999.4
$ws.Range("L67").Value = 9999.4
$ws.Range("N67").Value = -11715.4
$ws.Range("H70").Value = 7856.3716
$ws.Range("J70").Value = 7943.6875
$ws.Range("L70").Value = 23831.0625
$ws.Range("N70").Value = -24371.0625
$ws.Range("H73").Value = 7856.3716
$ws.Range("J73").Value = 7943.6875
$ws.Range("L73").Value = 23831.0625
$ws.Range("N73").Value = -25703.0625
$ws.Range("H74").Value = 6909.6665
$ws.Range("I74").Value = 5487.25
$ws.Range("J74").Value = 7620.875
$ws.Range("K74").Value = 5487.25
$ws.Range("L74").Value = 7620.875
$ws.Range("M74").Value = -4551.25
$ws.Range("N74").Value = -9492.875
$ws.Range("H77").Value = 6909.6665
$ws.Range("I77").Value = 5487.25
$ws.Range("J77").Value = 7620.875
$ws.Range("K77").Value = 27436.25
$ws.Range("L77").Value = 38104.375
$ws.Range("M77").Value = -22756.25
$ws.Range("N77").Value = -47464.375
$ws.Range("H86").Value = 6181.8
$ws.Range("J86").Value = 6740.5835
$ws.Range("L86").Value = 6740.5835
$ws.Range("N86").Value = -8986.583500000001
$ws.Range("H89").Value = 6181.8
$ws.Range("J89").Value = 6740.5835
$ws.Range("L89").Value = 33702.9175
$ws.Range("N89").Value = -44934.9175
$ws.Range("H101").Value = 27778172
$ws.Range("I101").Value = 41666976
$ws.Range("J101").Value = 566.6667
$ws.Range("K101").Value = 125000928
$ws.Range("L101").Value = 1700.0001
$ws.Range("M101").Value = -124999306
$ws.Range("N101").Value = -4944.0001
$ws.Range("H125").Value = 8549800
$ws.Range("I125").Value = 1405.3077
$ws.Range("J125").Value = 12823998
$ws.Range("K125").Value = 12647.7693
$ws.Range("L125").Value = 115415982
$ws.Range("M125").Value = -10187.7693
$ws.Range("N125").Value = -115420902
$ws.Range("H138").Value = 2952.3618
$ws.Range("J138").Value = 3794.75
$ws.Range("L138").Value = 11384.25
$ws.Range("N138").Value = -21664.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7541.921
$ws.Range("I32").Value = 4126.0635
$ws.Range("K32").Value = 4126.0635
$ws.Range("M32").Value = -3839.0635
$ws.Range("H61").Value = 6942
$ws.Range("I61").Value = 7498.4287
$ws.Range("J61").Value = 4994.5
$ws.Range("K61").Value = 7498.4287
$ws.Range("L61").Value = 4994.5
$ws.Range("M61").Value = -7286.4287
$ws.Range("N61").Value = -5418.5
$ws.Range("H63").Value = 4636.8335
$ws.Range("I63").Value = 2525.3076
$ws.Range("J63").Value = 7132.273
$ws.Range("K63").Value = 2525.3076
$ws.Range("L63").Value = 7132.273
$ws.Range("M63").Value = -1839.3076
$ws.Range("N63").Value = -8504.273000000001
$ws.Range("H66").Value = 4636.8335
$ws.Range("I66").Value = 2525.3076
$ws.Range("J66").Value = 7132.273
$ws.Range("K66").Value = 12626.538
$ws.Range("L66").Value = 35661.365
$ws.Range("M66").Value = -9194.538
$ws.Range("N66").Value = -42525.365
$ws.Range("H102").Value = 4905424.5
$ws.Range("I102").Value = 6412957
$ws.Range("K102").Value = 6412957
$ws.Range("M102").Value = -6411335
$ws.Range("H122").Value = 511321.78
$ws.Range("I122").Value = 2029.7
$ws.Range("K122").Value = 6089.1
$ws.Range("M122").Value = -3639.1
$ws.Range("H132").Value = 3093.7104
$ws.Range("I132").Value = 2386.7273
$ws.Range("K132").Value = 7160.1819
$ws.Range("M132").Value = -4630.1819
$ws.Range("H136").Value = 6942
$ws.Range("I136").Value = 7498.4287
$ws.Range("J136").Value = 4994.5
$ws.Range("K136").Value = 22495.2861
$ws.Range("L136").Value = 14983.5
$ws.Range("M136").Value = -19945.2861
$ws.Range("N136").Value = -20083.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4168189.5
$ws.Range("I86").Value = 6668283
$ws.Range("J86").Value = 1367.5555
$ws.Range("K86").Value = 6668283
$ws.Range("L86").Value = 1367.5555
$ws.Range("M86").Value = -6667160
$ws.Range("N86").Value = -3613.5555
$ws.Range("H89").Value = 4168189.5
$ws.Range("I89").Value = 6668283
$ws.Range("J89").Value = 1367.5555
$ws.Range("K89").Value = 33341415
$ws.Range("L89").Value = 6837.7775
$ws.Range("M89").Value = -33335799
$ws.Range("N89").Value = -18069.7775
$ws.Range("H134").Value = 3291.311
$ws.Range("I134").Value = 1583.3103
$ws.Range("K134").Value = 4749.9309
$ws.Range("M134").Value = -2214.9309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 39999.5
$ws.Range("J28").Value = 39999.5
$ws.Range("L28").Value = 39999.5
$ws.Range("N28").Value = -40489.5
$ws.Range("H31").Value = 32405.871
$ws.Range("I31").Value = 2527.0454
$ws.Range("J31").Value = 105443
$ws.Range("K31").Value = 2527.0454
$ws.Range("L31").Value = 105443
$ws.Range("M31").Value = -2232.0454
$ws.Range("N31").Value = -106033
$ws.Range("H34").Value = 32405.871
$ws.Range("I34").Value = 2527.0454
$ws.Range("J34").Value = 105443
$ws.Range("K34").Value = 2527.0454
$ws.Range("L34").Value = 105443
$ws.Range("M34").Value = -2325.0454
$ws.Range("N34").Value = -105847

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10320844
$ws.Range("I4").Value = 11297928
$ws.Range("K4").Value = 33893784
$ws.Range("M4").Value = -33893672
$ws.Range("H122").Value = 1267.4286
$ws.Range("I122").Value = 1267.4286
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11406.8574
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8956.857399999999
$ws.Range("H132").Value = 1800.1
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 1811.2222
$ws.Range("K132").Value = 15300
$ws.Range("L132").Value = 16300.9998
$ws.Range("M132").Value = -12770
$ws.Range("N132").Value = -21360.9998
$ws.Range("H133").Value = 2010
$ws.Range("I133").Value = 2010
$ws.Range("K133").Value = 6030
$ws.Range("M133").Value = -970
$ws.Range("H138").Value = 3571.5557
$ws.Range("J138").Value = 3700
$ws.Range("L138").Value = 11100
$ws.Range("N138").Value = -21380
$ws.Range("H139").Value = 1070
$ws.Range("I139").Value = 1070
$ws.Range("K139").Value = 3210
$ws.Range("M139").Value = 1930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 358266.7
$ws.Range("I122").Value = 525631.75
$ws.Range("J122").Value = 2615.875
$ws.Range("K122").Value = 1576895.25
$ws.Range("L122").Value = 7847.625
$ws.Range("M122").Value = -1574445.25
$ws.Range("N122").Value = -12747.625
$ws.Range("H132").Value = 3349.7856
$ws.Range("I132").Value = 3035.2173
$ws.Range("K132").Value = 9105.651899999999
$ws.Range("M132").Value = -6575.651899999999
$ws.Range("H139").Value = 94215.375
$ws.Range("J139").Value = 94215.375
$ws.Range("L139").Value = 94215.375
$ws.Range("N139").Value = -104495.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1346.4
$ws.Range("I16").Value = 1310.3334
$ws.Range("J16").Value = 1400.5
$ws.Range("K16").Value = 1310.3334
$ws.Range("L16").Value = 1400.5
$ws.Range("M16").Value = -1140.3334
$ws.Range("N16").Value = -1740.5
$ws.Range("H40").Value = 9331.166999999999
$ws.Range("I40").Value = 5495
$ws.Range("K40").Value = 5495
$ws.Range("M40").Value = -5359
$ws.Range("H46").Value = 6765.8335
$ws.Range("I46").Value = 2799.5
$ws.Range("K46").Value = 2799.5
$ws.Range("M46").Value = -2611.5
$ws.Range("H61").Value = 15876031
$ws.Range("I61").Value = 22222956
$ws.Range("K61").Value = 22222956
$ws.Range("M61").Value = -22222754
$ws.Range("H113").Value = 15876031
$ws.Range("I113").Value = 22222956
$ws.Range("K113").Value = 22222956
$ws.Range("M113").Value = -22220786

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2457.75
$ws.Range("I126").Value = 2687.125
$ws.Range("K126").Value = 8061.375
$ws.Range("M126").Value = -5591.375
$ws.Range("H132").Value = 31285818
$ws.Range("I132").Value = 37045860
$ws.Range("J132").Value = 181577.6
$ws.Range("K132").Value = 111137580
$ws.Range("L132").Value = 544732.8
$ws.Range("M132").Value = -111135050
$ws.Range("N132").Value = -549792.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N122").ClearContents()
